# Update the Min/Max threshold values on Sheet1 to match the new dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# alpha_distance_range (row 2): Min 3.8 -> 5, Max 11.5 -> 11.7
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 11.7

# beta_distance_range (row 3): Min 4.6 -> 5, Max 10.5 -> 10.4
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 10.4

# ratio_threshold_range (row 4): Min 0.7 -> 0.8
$ws.Range("B4").Value = 0.8

# pie_threshold_range (row 5): Max 20 -> 22
$ws.Range("C5").Value = 22
